$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 with new values (per diff) ---
# (row 8 = "extr1" ... row 15 = "extr8" in column B)

# Row 8 (extr1): C/D/E change
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Row 9 (extr2): only C changes
$ws.Cells.Item(9, 3).Value = 16

# Row 10 (extr3): C/D/E change
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11 (extr4): C/D change
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9

# Row 12 (extr5): C/E change
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 5).Value = $false

# Row 13 (extr6): D changes only
$ws.Cells.Item(13, 4).Value = 8

# Row 14 (extr7): C/D change
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11

# Row 15 (extr8): C/D change
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11

# --- Add two new rows (16, 17) for line7 / line8 ---
# Column A on the existing rows uses a bold/bordered/centered style (s="1");
# copy that formatting onto the new A16/A17 cells instead of rebuilding it
# property-by-property (which would mint stray unused style entries).

# Row 16
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "line7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

# Row 17
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "line8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true
